$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Before:
#   1) Earnings of a person depends with the age of the person
#   2) Earnings of a person depends on how many hours they work
#   3) Earnings of a person depends on the education level of the person
#   4) (empty paragraph carrying the hidden "_GoBack" bookmark)
#
# After:
#   1) There is a difference between average earnings of people who are
#      having 2 or few kids with respect to the people who are having more
#      than 2 kids.                           (now carries "_GoBack")
#   2) Earnings of a person depends with the age of the person
#   3) Earnings of a person depends on how many hours they work
#   4) Earnings of a person depends on the education level of the person
#      (text split into two runs, with the relocated bookmarkEnd between
#       them, and the paragraph turned into a proper numbered ListParagraph)
#   5) (new, empty trailing ListParagraph, same as the old paragraph 4 minus
#      the bookmark)
# ---------------------------------------------------------------------------

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'
$pPrList = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + $rPr + '</w:pPr>'
$pPrPlain = '<w:pPr><w:pStyle w:val="ListParagraph"/>' + $rPr + '</w:pPr>'

function New-PackageXml([string]$bodyXml) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyXml
</w:body>
</w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@
}

function Get-ParagraphAtStart($doc, $startPos) {
    foreach ($para in $doc.Paragraphs) {
        if ($para.Range.Start -eq $startPos) {
            return $para
        }
    }
    return $null
}

# --- Step 1: insert a brand-new bullet right before "...age of the person" -
$find1 = $d.Content
$null = $find1.Find.Execute("Earnings of a person depends with the age of the person", `
        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertStart = $find1.Start
$null = $d.Range($insertStart, $insertStart).InsertParagraphBefore()

# The paragraph that now occupies [insertStart, insertStart+1) is the new,
# empty paragraph; replace its whole contents (pPr + runs + bookmark) with
# raw OOXML so the two sentences remain two separate runs and the bookmark
# sits exactly at the start of the paragraph.
$newPara = Get-ParagraphAtStart $d $insertStart
$newParaRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

$firstBulletXml = "<w:p>$pPrList" + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    "<w:r>$rPr" + '<w:t xml:space="preserve">There is a </w:t></w:r>' + `
    "<w:r>$rPr" + '<w:t>difference between average earnings of people who are having 2 or few kids with respect to the people who are having more than 2 kids.</w:t></w:r>' + `
    '</w:p>'

$null = $newParaRange.InsertXML((New-PackageXml $firstBulletXml))

# --- Step 2: rebuild the "education level" bullet (split text, relocate the
#     bookmark end) together with the blank paragraph right after it (which
#     used to hold bookmarkStart/bookmarkEnd for "_GoBack"). ----------------
$find2 = $d.Content
$null = $find2.Find.Execute("Earnings of a person depends on the education level of the person", `
        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$eduStart = $find2.Start

$eduPara = Get-ParagraphAtStart $d $eduStart
$blankPara = Get-ParagraphAtStart $d $eduPara.Range.End

$replaceRange = $d.Range($eduPara.Range.Start, $blankPara.Range.End)

$eduBlockXml = "<w:p>$pPrList" + `
    "<w:r>$rPr" + '<w:t xml:space="preserve">Earnings of a person depends on the education level of the </w:t></w:r>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    "<w:r>$rPr" + '<w:t>person</w:t></w:r>' + `
    '</w:p>' + `
    "<w:p>$pPrPlain</w:p>"

$null = $replaceRange.InsertXML((New-PackageXml $eduBlockXml))
